# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Commit: Updated cryptos list on Tue Feb 27 09:57:41 UTC 2024 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.636.02'
$ws.Range("E2").Value = '  +11.15%  '
$ws.Range("D3").Value = '3.257.11'
$ws.Range("E3").Value = '  +7.11%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''399.12'
$ws.Range("E5").Value = '  +3.80%  '
$ws.Range("D6").Value = '''111.20'
$ws.Range("E6").Value = '  +9.55%  '
$ws.Range("E7").Value = '  +5.64%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '''0.622'
$ws.Range("E9").Value = '  +8.39%  '
$ws.Range("D10").Value = '''39.42'
$ws.Range("E10").Value = '  +8.52%  '
$ws.Range("D11").Value = '''0.0959'
$ws.Range("E11").Value = '  +13.72%  '
$ws.Range("E12").Value = '  +2.73%  '
$ws.Range("D13").Value = '3.764.37'
$ws.Range("E13").Value = '  +6.56%  '
$ws.Range("D14").Value = '''19.26'
$ws.Range("E14").Value = '  +5.99%  '
$ws.Range("D15").Value = '''8.14'
$ws.Range("E15").Value = '  +6.80%  '
$ws.Range("D16").Value = '3.252.06'
$ws.Range("E16").Value = '  +6.92%  '
$ws.Range("E17").Value = '  +7.78%  '
$ws.Range("D18").Value = '''11.09'
$ws.Range("E18").Value = '  +4.49%  '
$ws.Range("D19").Value = '56.511.86'
$ws.Range("E19").Value = '  +10.78%  '
$ws.Range("E20").Value = '  +4.99%  '
$ws.Range("E21").Value = '  +10.50%  '
$ws.Range("D22").Value = '''13.09'
$ws.Range("E22").Value = '  +7.86%  '
$ws.Range("D23").Value = '''300.38'
$ws.Range("E23").Value = '  +14.35%  '
$ws.Range("D24").Value = '''75.62'
$ws.Range("E24").Value = '  +9.18%  '
$ws.Range("E25").Value = '  +4.92%  '
$ws.Range("D26").Value = '''8.21'
$ws.Range("E26").Value = '  +4.92%  '
$ws.Range("D27").Value = '''28.46'
$ws.Range("E27").Value = '  +6.33%  '
$ws.Range("E28").Value = '  +4.18%  '
$ws.Range("E29").Value = '  +3.55%  '
$ws.Range("D30").Value = '''0.172'
$ws.Range("E30").Value = '  +6.56%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  +7.46%  '
$ws.Range("D33").Value = '''11.15'
$ws.Range("E33").Value = '  +8.33%  '
$ws.Range("D34").Value = '''36.91'
$ws.Range("E34").Value = '  +5.79%  '
$ws.Range("E35").Value = '  +5.14%  '
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("D37").Value = '''51.50'
$ws.Range("E37").Value = '  +3.32%  '
$ws.Range("E38").Value = '  +6.65%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''1.00'
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''3.09'
$ws.Range("E40").Value = '  +26.83%  '
$ws.Range("D41").Value = '''17.65'
$ws.Range("E41").Value = '  +8.87%  '
$ws.Range("D42").Value = '''134.65'
$ws.Range("E42").Value = '  +3.84%  '
$ws.Range("D43").Value = '''1.94'
$ws.Range("E43").Value = '  +6.83%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.120'
$ws.Range("E44").Value = '  +5.35%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''4.00'
$ws.Range("E45").Value = '  +7.82%  '
$ws.Range("E46").Value = '  -0.80%  '
$ws.Range("D47").Value = '''22.38'
$ws.Range("E47").Value = '  +4.51%  '
$ws.Range("D48").Value = '''2.19'
$ws.Range("E48").Value = '  +58.94%  '
$ws.Range("D49").Value = '2.145.40'
$ws.Range("E49").Value = '  +5.12%  '
$ws.Range("D50").Value = '''2.08'
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("D51").Value = '''2.42'
$ws.Range("E51").Value = '  -2.83%  '
